$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.49967933333333
$ws.Range("H2").Value = 58.499038
$ws.Range("I2").Value = 0.01453409262904611
$ws.Range("J2").Value = 0.01453409262904611
$ws.Range("M2").Value = 0.1808983333333333
$ws.Range("N2").Value = 0.5426949999999999
$ws.Range("O2").Value = 0.09546831801815302
$ws.Range("P2").Value = 0.09546831801815302
$ws.Range("Q2").Value = 3.527459491934444
$ws.Range("R2").Value = 31.74713542740999
$ws.Range("S2").Value = 0.001387545377215068
$ws.Range("T2").Value = 0.001387545377215068
$ws.Range("G3").Value = 19.49967933333333
$ws.Range("H3").Value = 58.499038
$ws.Range("I3").Value = 0.01453409262904611
$ws.Range("J3").Value = 0.01453409262904611
$ws.Range("M3").Value = 1.572737
$ws.Range("N3").Value = 4.718211
$ws.Range("O3").Value = 0.8300051930177132
$ws.Range("P3").Value = 0.8300051930177132
$ws.Range("Q3").Value = 30.66786717566867
$ws.Range("R3").Value = 276.010804581018
$ws.Range("S3").Value = 0.01206337235790874
$ws.Range("T3").Value = 0.01206337235790874
$ws.Range("G4").Value = 19.49967933333333
$ws.Range("H4").Value = 58.499038
$ws.Range("I4").Value = 0.01453409262904611
$ws.Range("J4").Value = 0.01453409262904611
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1412166666666667
$ws.Range("N4").Value = 0.42365
$ws.Range("O4").Value = 0.07452648896413371
$ws.Range("P4").Value = 0.07452648896413369
$ws.Range("Q4").Value = 2.753679716522222
$ws.Range("R4").Value = 24.7831174487
$ws.Range("S4").Value = 0.001083174893922302
$ws.Range("T4").Value = 0.001083174893922302
$ws.Range("I5").Value = 0.2821439310161206
$ws.Range("J5").Value = 0.2821439310161206
$ws.Range("M5").Value = 0.1808983333333333
$ws.Range("N5").Value = 0.5426949999999999
$ws.Range("O5").Value = 0.09546831801815302
$ws.Range("P5").Value = 0.09546831801815302
$ws.Range("Q5").Value = 68.47701559060664
$ws.Range("R5").Value = 616.2931403154598
$ws.Range("S5").Value = 0.02693580653313882
$ws.Range("T5").Value = 0.02693580653313883
$ws.Range("I6").Value = 0.2821439310161206
$ws.Range("J6").Value = 0.2821439310161206
$ws.Range("M6").Value = 1.572737
$ws.Range("N6").Value = 4.718211
$ws.Range("O6").Value = 0.8300051930177132
$ws.Range("P6").Value = 0.8300051930177132
$ws.Range("Q6").Value = 595.3417816762119
$ws.Range("R6").Value = 5358.076035085908
$ws.Range("S6").Value = 0.2341809279218115
$ws.Range("T6").Value = 0.2341809279218116
$ws.Range("I7").Value = 0.2821439310161206
$ws.Range("J7").Value = 0.2821439310161206
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1412166666666667
$ws.Range("N7").Value = 0.42365
$ws.Range("O7").Value = 0.07452648896413371
$ws.Range("P7").Value = 0.07452648896413369
$ws.Range("Q7").Value = 53.45597002913333
$ws.Range("R7").Value = 481.1037302622
$ws.Range("S7").Value = 0.02102719656117021
$ws.Range("T7").Value = 0.02102719656117021
$ws.Range("G8").Value = 481.5587156666667
$ws.Range("H8").Value = 1444.676147
$ws.Range("I8").Value = 0.3589299526510408
$ws.Range("J8").Value = 0.3589299526510408
$ws.Range("M8").Value = 0.1808983333333333
$ws.Range("N8").Value = 0.5426949999999999
$ws.Range("O8").Value = 0.09546831801815302
$ws.Range("P8").Value = 0.09546831801815302
$ws.Range("Q8").Value = 87.11316906624054
$ws.Range("R8").Value = 784.018521596165
$ws.Range("S8").Value = 0.03426643886593016
$ws.Range("T8").Value = 0.03426643886593016
$ws.Range("G9").Value = 481.5587156666667
$ws.Range("H9").Value = 1444.676147
$ws.Range("I9").Value = 0.3589299526510408
$ws.Range("J9").Value = 0.3589299526510408
$ws.Range("M9").Value = 1.572737
$ws.Range("N9").Value = 4.718211
$ws.Range("O9").Value = 0.8300051930177132
$ws.Range("P9").Value = 0.8300051930177132
$ws.Range("Q9").Value = 757.3652098014464
$ws.Range("R9").Value = 6816.286888213018
$ws.Range("S9").Value = 0.2979137246299658
$ws.Range("T9").Value = 0.2979137246299658
$ws.Range("G10").Value = 481.5587156666667
$ws.Range("H10").Value = 1444.676147
$ws.Range("I10").Value = 0.3589299526510408
$ws.Range("J10").Value = 0.3589299526510408
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1412166666666667
$ws.Range("N10").Value = 0.42365
$ws.Range("O10").Value = 0.07452648896413371
$ws.Range("P10").Value = 0.07452648896413369
$ws.Range("Q10").Value = 68.00411663072779
$ws.Range("R10").Value = 612.0370496765501
$ws.Range("S10").Value = 0.02674978915514482
$ws.Range("T10").Value = 0.02674978915514482
$ws.Range("G11").Value = 10.909999
$ws.Range("H11").Value = 32.729997
$ws.Range("I11").Value = 0.008131771468556478
$ws.Range("J11").Value = 0.008131771468556478
$ws.Range("M11").Value = 0.1808983333333333
$ws.Range("N11").Value = 0.5426949999999999
$ws.Range("O11").Value = 0.09546831801815302
$ws.Range("P11").Value = 0.09546831801815302
$ws.Range("Q11").Value = 1.973600635768333
$ws.Range("R11").Value = 17.76240572191499
$ws.Range("S11").Value = 0.000776326544611093
$ws.Range("T11").Value = 0.000776326544611093
$ws.Range("G12").Value = 10.909999
$ws.Range("H12").Value = 32.729997
$ws.Range("I12").Value = 0.008131771468556478
$ws.Range("J12").Value = 0.008131771468556478
$ws.Range("M12").Value = 1.572737
$ws.Range("N12").Value = 4.718211
$ws.Range("O12").Value = 0.8300051930177132
$ws.Range("P12").Value = 0.8300051930177132
$ws.Range("Q12").Value = 17.158559097263
$ws.Range("R12").Value = 154.427031875367
$ws.Range("S12").Value = 0.006749412547335153
$ws.Range("T12").Value = 0.006749412547335153
$ws.Range("G13").Value = 10.909999
$ws.Range("H13").Value = 32.729997
$ws.Range("I13").Value = 0.008131771468556478
$ws.Range("J13").Value = 0.008131771468556478
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1412166666666667
$ws.Range("N13").Value = 0.42365
$ws.Range("O13").Value = 0.07452648896413371
$ws.Range("P13").Value = 0.07452648896413369
$ws.Range("Q13").Value = 1.540673692116667
$ws.Range("R13").Value = 13.86606322905
$ws.Range("S13").Value = 0.0006060323766102317
$ws.Range("T13").Value = 0.0006060323766102316
$ws.Range("G14").Value = 98.48487833333333
$ws.Range("H14").Value = 295.454635
$ws.Range("I14").Value = 0.07340573759129181
$ws.Range("J14").Value = 0.07340573759129182
$ws.Range("M14").Value = 0.1808983333333333
$ws.Range("N14").Value = 0.5426949999999999
$ws.Range("O14").Value = 0.09546831801815302
$ws.Range("P14").Value = 0.09546831801815302
$ws.Range("Q14").Value = 17.81575034903611
$ws.Range("R14").Value = 160.341753141325
$ws.Range("S14").Value = 0.007007922300722536
$ws.Range("T14").Value = 0.007007922300722537
$ws.Range("G15").Value = 98.48487833333333
$ws.Range("H15").Value = 295.454635
$ws.Range("I15").Value = 0.07340573759129181
$ws.Range("J15").Value = 0.07340573759129182
$ws.Range("M15").Value = 1.572737
$ws.Range("N15").Value = 4.718211
$ws.Range("O15").Value = 0.8300051930177132
$ws.Range("P15").Value = 0.8300051930177132
$ws.Range("Q15").Value = 154.8908120953317
$ws.Range("R15").Value = 1394.017308857985
$ws.Range("S15").Value = 0.06092714339806776
$ws.Range("T15").Value = 0.06092714339806778
$ws.Range("G16").Value = 98.48487833333333
$ws.Range("H16").Value = 295.454635
$ws.Range("I16").Value = 0.07340573759129181
$ws.Range("J16").Value = 0.07340573759129182
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.1412166666666667
$ws.Range("N16").Value = 0.42365
$ws.Range("O16").Value = 0.07452648896413371
$ws.Range("P16").Value = 0.07452648896413369
$ws.Range("Q16").Value = 13.90770623530556
$ws.Range("R16").Value = 125.16935611775
$ws.Range("S16").Value = 0.005470671892501503
$ws.Range("T16").Value = 0.005470671892501503
$ws.Range("G17").Value = 352.659012
$ws.Range("H17").Value = 1057.977036
$ws.Range("I17").Value = 0.2628545146439442
$ws.Range("J17").Value = 0.2628545146439442
$ws.Range("M17").Value = 0.1808983333333333
$ws.Range("N17").Value = 0.5426949999999999
$ws.Range("O17").Value = 0.09546831801815302
$ws.Range("P17").Value = 0.09546831801815302
$ws.Range("Q17").Value = 63.79542750577999
$ws.Range("R17").Value = 574.1588475520199
$ws.Range("S17").Value = 0.02509427839653533
$ws.Range("T17").Value = 0.02509427839653533
$ws.Range("G18").Value = 352.659012
$ws.Range("H18").Value = 1057.977036
$ws.Range("I18").Value = 0.2628545146439442
$ws.Range("J18").Value = 0.2628545146439442
$ws.Range("M18").Value = 1.572737
$ws.Range("N18").Value = 4.718211
$ws.Range("O18").Value = 0.8300051930177132
$ws.Range("P18").Value = 0.8300051930177132
$ws.Range("Q18").Value = 554.639876555844
$ws.Range("R18").Value = 4991.758889002596
$ws.Range("S18").Value = 0.2181706121626243
$ws.Range("T18").Value = 0.2181706121626243
$ws.Range("G19").Value = 352.659012
$ws.Range("H19").Value = 1057.977036
$ws.Range("I19").Value = 0.2628545146439442
$ws.Range("J19").Value = 0.2628545146439442
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.1412166666666667
$ws.Range("N19").Value = 0.42365
$ws.Range("O19").Value = 0.07452648896413371
$ws.Range("P19").Value = 0.07452648896413369
$ws.Range("Q19").Value = 49.80133014460001
$ws.Range("R19").Value = 448.2119713014
$ws.Range("S19").Value = 0.01958962408478463
$ws.Range("T19").Value = 0.01958962408478463
